# Update imputed values in the KNN result data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.6
$ws.Range("A6").Value = -22.056
$ws.Range("A7").Value = -21.286
$ws.Range("B7").Value = 6.312
$ws.Range("B12").Value = 5.401999999999999
$ws.Range("E13").Value = 16.561
$ws.Range("E14").Value = 17.048
$ws.Range("B15").Value = 5.081
$ws.Range("A16").Value = -21.397
$ws.Range("E16").Value = 16.784
$ws.Range("E19").Value = 16.694
$ws.Range("A20").Value = -21.807
$ws.Range("B20").Value = 5.532999999999999
$ws.Range("B21").Value = 8.564
$ws.Range("B22").Value = 7.034000000000001
$ws.Range("E22").Value = 16.493
$ws.Range("B23").Value = 7.007
$ws.Range("A28").Value = -21.878
$ws.Range("A29").Value = -21.664
$ws.Range("B29").Value = 5.965999999999999
$ws.Range("A32").Value = -21.705
$ws.Range("B34").Value = 8.059000000000001
$ws.Range("E36").Value = 16.748
$ws.Range("A40").Value = -20.312
$ws.Range("B42").Value = 7.238
$ws.Range("B43").Value = 5.529000000000001
$ws.Range("B44").Value = 5.08
$ws.Range("B45").Value = 5.315
$ws.Range("A46").Value = -21.017
$ws.Range("B46").Value = 6.789
$ws.Range("E46").Value = 16.695
$ws.Range("B50").Value = 5.584000000000001
$ws.Range("E50").Value = 16.514
$ws.Range("A51").Value = -21.079
$ws.Range("B51").Value = 7.498
$ws.Range("A52").Value = -21.261
$ws.Range("A57").Value = -22.354
$ws.Range("A59").Value = -22.425
$ws.Range("A62").Value = -21.676
$ws.Range("A66").Value = -21.504
$ws.Range("B66").Value = 5.709
$ws.Range("B67").Value = 5.430999999999999
$ws.Range("A73").Value = -20.53
$ws.Range("A74").Value = -20.947
$ws.Range("B79").Value = 5.680999999999999
$ws.Range("B84").Value = 5.781000000000001
$ws.Range("A92").Value = -21.299
$ws.Range("B92").Value = 5.681000000000001
$ws.Range("E95").Value = 17.159
$ws.Range("B97").Value = 5.677
$ws.Range("E97").Value = 16.922
$ws.Range("A100").Value = -21.481
